# Settings file added and read.
# Applies the "Data" sheet schema expansion (Yahoo Finance style columns)
# and adds an instructional note block on the "Aktier" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Aktier": update labels, add hint text + instructional note box
# ---------------------------------------------------------------------
$aktier = $wb.Worksheets.Item("Aktier")

$aktier.Range("A1").Value = "Namn"
$aktier.Range("B1").Value = "Förkortning"
$aktier.Range("A2").Value = "Volvo B"
$aktier.Range("B2").Value = "VOLV-B.ST"
$aktier.Range("A3").Value = "HM"
$aktier.Range("B3").Value = "HM-B.ST"

$aktier.Range("F2").Value = "Hämta förkortningen från https://finance.yahoo.com/"

$note = $aktier.Range("F3:J5")
$note.Merge()
$note.Value = "Där ser du också vilken data som kan hämtas. Endast den som har kolumner i Data arket sparas men modifikationer till programmet kan göras enkelt. "
$note.HorizontalAlignment = -4108
$note.WrapText = $true

$aktier.Range("F8").Select()

# ---------------------------------------------------------------------
# Sheet "Data": expand table with the full set of Yahoo Finance fields
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")
$lo = $data.ListObjects.Item(1)

$lo.Resize($data.Range("A1:H4"))

$data.Range("A1").Value = "Aktie beteckning"
$data.Range("B1").Value = "TradeDate "
$data.Range("C1").Value = "EarningsShare "
$data.Range("D1").Value = "DaysLow "
$data.Range("E1").Value = "DaysHigh "
$data.Range("F1").Value = "LastTradePriceOnly "
$data.Range("G1").Value = "Open "
$data.Range("H1").Value = "DateStamp"

$headerRow = $data.Range("A1:H1")
$headerRow.Font.Bold = $true
$headerRow.Borders.Item(3).LineStyle = 1
$headerRow.Borders.Item(3).Weight = 2
$headerRow.Borders.Item(4).LineStyle = 1
$headerRow.Borders.Item(4).Weight = 2

$data.Range("H1:H3").NumberFormat = "yyyy/mm/dd\ hh:mm;@"

# Row 2 - HM-B.ST
$data.Range("A2").Value = "HM-B.ST"
$data.Range("B2").Value = ""
$data.Range("C2").Value = ""
$data.Range("D2").Value = "259.60"
$data.Range("E2").Value = "265.90"
$data.Range("F2").Value = "265.90"
$data.Range("G2").Value = "260.00"
$data.Range("H2").Value = 42691.6109973958

# Row 3 - VOLV-B.ST
$data.Range("A3").Value = "VOLV-B.ST"
$data.Range("B3").Value = ""
$data.Range("C3").Value = "7.17"
$data.Range("D3").Value = "95.55"
$data.Range("E3").Value = "97.25"
$data.Range("F3").Value = "96.10"
$data.Range("G3").Value = "96.10"
$data.Range("H3").Value = 42691.6109959722

$data.Range("A2:A2").EntireRow.Select()
